# Updates the cryptocurrency price/volume table in the worksheet to match
# the latest scrape. Most rows only get their "Price" (column D) and
# "Volume(1h)" (column E) values refreshed; a few rows additionally shift
# identity (Coin name + Link) because the underlying ranking reordered
# (rows 42-44: ApeXProtocol / dogwifhat / Kaspa).
#
# Price values are written with a leading "'" where the text looks like a
# pure number (e.g. "579.12") so Excel keeps storing them as text, matching
# the original inlineStr/text representation of that column instead of
# silently converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.154.34"
$ws.Cells.Item(2, 5).Value = "  -0.74%  "
$ws.Cells.Item(3, 4).Value = "3.579.84"
$ws.Cells.Item(3, 5).Value = "  -1.46%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "'579.12"
$ws.Cells.Item(5, 5).Value = "  -2.14%  "
$ws.Cells.Item(6, 4).Value = "'186.92"
$ws.Cells.Item(6, 5).Value = "  -4.13%  "
$ws.Cells.Item(7, 4).Value = "3.576.86"
$ws.Cells.Item(7, 5).Value = "  -1.33%  "
$ws.Cells.Item(8, 5).Value = "  -3.63%  "
$ws.Cells.Item(9, 4).Value = "'0.999"
$ws.Cells.Item(9, 5).Value = "  -0.05%  "
$ws.Cells.Item(10, 5).Value = "  -0.87%  "
$ws.Cells.Item(11, 4).Value = "'0.653"
$ws.Cells.Item(11, 5).Value = "  -3.93%  "
$ws.Cells.Item(12, 4).Value = "'55.17"
$ws.Cells.Item(12, 5).Value = "  -4.92%  "
$ws.Cells.Item(13, 5).Value = "  -2.27%  "
$ws.Cells.Item(14, 4).Value = "'9.55"
$ws.Cells.Item(14, 5).Value = "  -4.01%  "
$ws.Cells.Item(15, 4).Value = "4.152.51"
$ws.Cells.Item(15, 5).Value = "  -1.54%  "
$ws.Cells.Item(16, 4).Value = "'19.66"
$ws.Cells.Item(16, 5).Value = "  -4.35%  "
$ws.Cells.Item(17, 4).Value = "3.567.93"
$ws.Cells.Item(17, 5).Value = "  -1.91%  "
$ws.Cells.Item(18, 4).Value = "70.038.17"
$ws.Cells.Item(18, 5).Value = "  -0.94%  "
$ws.Cells.Item(19, 4).Value = "'12.59"
$ws.Cells.Item(19, 5).Value = "  -1.18%  "
$ws.Cells.Item(20, 5).Value = "  -1.03%  "
$ws.Cells.Item(21, 5).Value = "  -3.00%  "
$ws.Cells.Item(22, 4).Value = "'493.24"
$ws.Cells.Item(22, 5).Value = "  +0.91%  "
$ws.Cells.Item(23, 4).Value = "'19.50"
$ws.Cells.Item(23, 5).Value = "  +0.89%  "
$ws.Cells.Item(24, 4).Value = "'4.96"
$ws.Cells.Item(24, 5).Value = "  -5.48%  "
$ws.Cells.Item(25, 4).Value = "'97.26"
$ws.Cells.Item(25, 5).Value = "  +6.48%  "
$ws.Cells.Item(26, 4).Value = "'4.39"
$ws.Cells.Item(26, 5).Value = "  -2.33%  "
$ws.Cells.Item(27, 4).Value = "'11.46"
$ws.Cells.Item(27, 5).Value = "  +0.34%  "
$ws.Cells.Item(28, 5).Value = "  -6.32%  "
$ws.Cells.Item(29, 5).Value = "  -2.40%  "
$ws.Cells.Item(30, 4).Value = "'7.70"
$ws.Cells.Item(30, 5).Value = "  -3.12%  "
$ws.Cells.Item(31, 4).Value = "'31.73"
$ws.Cells.Item(31, 5).Value = "  -3.29%  "
$ws.Cells.Item(32, 4).Value = "'12.14"
$ws.Cells.Item(32, 5).Value = "  -1.25%  "
$ws.Cells.Item(33, 4).Value = "'65.89"
$ws.Cells.Item(33, 5).Value = "  -0.63%  "
$ws.Cells.Item(34, 5).Value = "  -5.64%  "
$ws.Cells.Item(35, 4).Value = "'573.14"
$ws.Cells.Item(35, 5).Value = "  -6.40%  "
$ws.Cells.Item(36, 4).Value = "'3.20"
$ws.Cells.Item(36, 5).Value = "  +12.70%  "
$ws.Cells.Item(37, 4).Value = "'39.08"
$ws.Cells.Item(37, 5).Value = "  -3.99%  "
$ws.Cells.Item(38, 4).Value = "'0.407"
$ws.Cells.Item(38, 5).Value = "  -1.14%  "
$ws.Cells.Item(39, 5).Value = "  +0.08%  "
$ws.Cells.Item(40, 4).Value = "0.0₃0796"
$ws.Cells.Item(40, 5).Value = "  -5.15%  "
$ws.Cells.Item(41, 4).Value = "'3.49"
$ws.Cells.Item(41, 5).Value = "  -2.40%  "
$ws.Cells.Item(42, 2).Value = "ApeXProtocol"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(42, 4).Value = "'3.72"
$ws.Cells.Item(42, 5).Value = "  +11.58%  "
$ws.Cells.Item(43, 2).Value = "dogwifhat"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(43, 4).Value = "'3.16"
$ws.Cells.Item(43, 5).Value = "  -1.76%  "
$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).Value = "'0.134"
$ws.Cells.Item(44, 5).Value = "  -10.00%  "
$ws.Cells.Item(45, 4).Value = "'3.08"
$ws.Cells.Item(45, 5).Value = "  -3.38%  "
$ws.Cells.Item(46, 5).Value = "  -0.87%  "
$ws.Cells.Item(47, 4).Value = "3.187.59"
$ws.Cells.Item(47, 5).Value = "  -4.02%  "
$ws.Cells.Item(48, 4).Value = "'9.51"
$ws.Cells.Item(48, 5).Value = "  -1.82%  "
$ws.Cells.Item(49, 4).Value = "'1.57"
$ws.Cells.Item(49, 5).Value = "  +31.18%  "
$ws.Cells.Item(50, 5).Value = "  -2.23%  "
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
